$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA14").Value = -1
$ws.Range("AB14").Value = 1
$ws.Range("AC14").Value = 0.925
$ws.Range("AD14").Value = -1
$ws.Range("B14").Value = 6876471
$ws.Range("E14").Value = "Portimonense"
$ws.Range("F14").Value = "Boavista"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 4
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = "A"
$ws.Range("L14").Value = 2.45
$ws.Range("M14").Value = 3.25
$ws.Range("N14").Value = 2.875
$ws.Range("O14").Value = 3
$ws.Range("P14").Value = 3.3
$ws.Range("Q14").Value = 2.3
$ws.Range("R14").Value = 0.25
$ws.Range("S14").Value = 1.85
$ws.Range("T14").Value = 2
$ws.Range("U14").Value = 2.25
$ws.Range("V14").Value = 1.925
$ws.Range("W14").Value = 1.925
$ws.Range("X14").Value = -1
$ws.Range("Z14").Value = 1.3
$ws.Range("AA15").Value = -0.5
$ws.Range("AB15").Value = 0.4625
$ws.Range("AC15").Value = -1
$ws.Range("AD15").Value = 0.9750000000000001
$ws.Range("B15").Value = 6876465
$ws.Range("E15").Value = "Benfica"
$ws.Range("F15").Value = "Estrela"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = "H"
$ws.Range("L15").Value = 1.125
$ws.Range("M15").Value = 8
$ws.Range("N15").Value = 19
$ws.Range("O15").Value = 1.142
$ws.Range("P15").Value = 7.5
$ws.Range("Q15").Value = 17
$ws.Range("R15").Value = -2.25
$ws.Range("S15").Value = 1.925
$ws.Range("T15").Value = 1.925
$ws.Range("U15").Value = 3.5
$ws.Range("V15").Value = 1.875
$ws.Range("W15").Value = 1.975
$ws.Range("X15").Value = 0.1419999999999999
$ws.Range("Z15").Value = -1
$ws.Range("AA94").Value = -0.5
$ws.Range("AB94").Value = 0.425
$ws.Range("B94").Value = 6876543
$ws.Range("E94").Value = "Vizela"
$ws.Range("F94").Value = "Famalicao"
$ws.Range("H94").Value = 0
$ws.Range("K94").Value = "D"
$ws.Range("L94").Value = 2.3
$ws.Range("N94").Value = 3.1
$ws.Range("O94").Value = 2.3
$ws.Range("P94").Value = 3.1
$ws.Range("Q94").Value = 3.1
$ws.Range("R94").Value = -0.25
$ws.Range("S94").Value = 2.05
$ws.Range("T94").Value = 1.85
$ws.Range("U94").Value = 2.25
$ws.Range("Y94").Value = 2.1
$ws.Range("Z94").Value = -1
$ws.Range("AA95").Value = -1
$ws.Range("AB95").Value = 0.875
$ws.Range("B95").Value = 6875464
$ws.Range("E95").Value = "Estrela"
$ws.Range("F95").Value = "Moreirense"
$ws.Range("H95").Value = 1
$ws.Range("K95").Value = "A"
$ws.Range("L95").Value = 2.5
$ws.Range("N95").Value = 2.875
$ws.Range("O95").Value = 2.625
$ws.Range("P95").Value = 3.3
$ws.Range("Q95").Value = 2.625
$ws.Range("R95").Value = 0
$ws.Range("S95").Value = 1.975
$ws.Range("T95").Value = 1.875
$ws.Range("U95").Value = 2.5
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 1.625
$ws.Range("AA151").Value = 1.025
$ws.Range("AB151").Value = -1
$ws.Range("AC151").Value = 0.925
$ws.Range("B151").Value = 6876586
$ws.Range("E151").Value = "Benfica"
$ws.Range("F151").Value = "Rio Ave"
$ws.Range("G151").Value = 4
$ws.Range("H151").Value = 1
$ws.Range("I151").Value = 1
$ws.Range("J151").Value = 1
$ws.Range("K151").Value = "H"
$ws.Range("L151").Value = 1.166
$ws.Range("M151").Value = 7.5
$ws.Range("N151").Value = 15
$ws.Range("O151").Value = 1.2
$ws.Range("P151").Value = 8
$ws.Range("Q151").Value = 10
$ws.Range("R151").Value = -2
$ws.Range("U151").Value = 3.25
$ws.Range("V151").Value = 1.925
$ws.Range("W151").Value = 1.925
$ws.Range("X151").Value = 0.2
$ws.Range("Z151").Value = -1
$ws.Range("AA152").Value = -1
$ws.Range("AB152").Value = 0.825
$ws.Range("AC152").Value = 1
$ws.Range("B152").Value = 6876591
$ws.Range("E152").Value = "Vizela"
$ws.Range("F152").Value = "Boavista"
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 4
$ws.Range("I152").Value = 0
$ws.Range("J152").Value = 2
$ws.Range("K152").Value = "A"
$ws.Range("L152").Value = 2.3
$ws.Range("M152").Value = 3.1
$ws.Range("N152").Value = 3.25
$ws.Range("O152").Value = 1.95
$ws.Range("P152").Value = 3.2
$ws.Range("Q152").Value = 4.2
$ws.Range("R152").Value = -0.5
$ws.Range("U152").Value = 2.25
$ws.Range("V152").Value = 2
$ws.Range("W152").Value = 1.85
$ws.Range("X152").Value = -1
$ws.Range("Z152").Value = 3.2
